# Generate Report for Handoff
# Refresh the localization-status report with a newly-generated handoff:
#  - new source file UUIDs (f5b2ccdf-... / ffffd6795c51-...)
#  - status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#  - new handoff timestamps, handback cleared back to the zero-date sentinel
#  - the now-redundant "handoff" columns (F/G) are dropped from the per-locale sheets

$wb = $excel.ActiveWorkbook

$oldUuid1 = "93bd694f-a1a8-4d8e-9745-eb97625df1f4"
$oldUuid2 = "a8eed73e-ebd3-4b54-8e2a-535cfbadb727"
$newUuid1 = "f5b2ccdf-10b2-4339-94e9-4b45c024f529"
$newUuid2 = "ffffd6795c51-6fb3-43a7-b7c7-a7226ad88f57"

$oldHash = "703667be1ac7faea0a185399db1e2c14e271649a"
$newHash = "ee3dcb8c376b760ac1fd1be7b860a5fba442b146"

$newMd1 = "$newUuid1.md"
$newMd2 = "$newUuid2.md"

$newXlfZh = "$newUuid1.$newHash.zh-cn.xlf"
$newXlfDe = "$newUuid1.$newHash.de-de.xlf"

$status = "Ready for handoff"
$newHandoffDate = "2016-03-23 17:16:29"
$newHandoffDatetimeZh = "2016-03-23 17:16:25"
$newHandbackDatetime = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $newMd1
$ws1.Range("B2").Value = $status
$ws1.Range("C2").Value = $status
$ws1.Range("D2").Value = $newHandoffDate

$ws1.Range("A3").Value = $newMd2
$ws1.Range("B3").Value = $status
$ws1.Range("C3").Value = $status
$ws1.Range("D3").Value = $newHandoffDate

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/520c068ea6794c12c29867ed60d34b6addaee20a/e2e/$newMd1", "", "", $newMd1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/520c068ea6794c12c29867ed60d34b6addaee20a/e2e/$newMd2", "", "", $newMd2) | Out-Null
$ws1.Range("A2").Style = "HyperLink"
$ws1.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newMd1
$ws2.Range("C2").Value = $status
$ws2.Range("D2").Value = $newXlfZh
$ws2.Range("E2").Value = $newHandoffDatetimeZh
$ws2.Range("H2").Value = $newHandbackDatetime

$ws2.Range("A3").Value = $newMd2
$ws2.Range("C3").Value = $status
$ws2.Range("D3").Value = $newXlfZh
$ws2.Range("E3").Value = $newHandoffDatetimeZh
$ws2.Range("H3").Value = $newHandbackDatetime

# The handoff-reference columns (F/G) are no longer needed - drop the cells.
$ws2.Range("F2:G3").Clear()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/520c068ea6794c12c29867ed60d34b6addaee20a/e2e/$newMd1", "", "", $newMd1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/800f9efe17cfd29e6dcdcb13d40586725209f9c6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZh", "", "", $newXlfZh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/520c068ea6794c12c29867ed60d34b6addaee20a/e2e/$newMd2", "", "", $newMd2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/800f9efe17cfd29e6dcdcb13d40586725209f9c6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZh", "", "", $newXlfZh) | Out-Null
$ws2.Range("A2").Style = "HyperLink"
$ws2.Range("D2").Style = "HyperLink"
$ws2.Range("A3").Style = "HyperLink"
$ws2.Range("D3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newMd1
$ws3.Range("C2").Value = $status
$ws3.Range("D2").Value = $newXlfDe
$ws3.Range("E2").Value = $newHandoffDate
$ws3.Range("H2").Value = $newHandbackDatetime

$ws3.Range("A3").Value = $newMd2
$ws3.Range("C3").Value = $status
$ws3.Range("D3").Value = $newXlfDe
$ws3.Range("E3").Value = $newHandoffDate
$ws3.Range("H3").Value = $newHandbackDatetime

# The handoff-reference columns (F/G) are no longer needed - drop the cells.
$ws3.Range("F2:G3").Clear()

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/520c068ea6794c12c29867ed60d34b6addaee20a/e2e/$newMd1", "", "", $newMd1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2b5f1a0a4045d7fcb53e47aca3356ae0ddb95cf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe", "", "", $newXlfDe) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/520c068ea6794c12c29867ed60d34b6addaee20a/e2e/$newMd2", "", "", $newMd2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d2b5f1a0a4045d7fcb53e47aca3356ae0ddb95cf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe", "", "", $newXlfDe) | Out-Null
$ws3.Range("A2").Style = "HyperLink"
$ws3.Range("D2").Style = "HyperLink"
$ws3.Range("A3").Style = "HyperLink"
$ws3.Range("D3").Style = "HyperLink"
